# "Change notes of quarterly series"
# Updates the footnote text in several rows of the "Fuentes y unidades" sheet,
# adjusts column widths for D:F, and updates the scroll/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the order in which new note text is first written matters, because it
# determines the order new entries are appended to the shared-strings table
# (xl/sharedStrings.xml), which must match the target workbook exactly.

# --- Rows 58 & 59: update note wording ---
$ws.Range("E58").Value = "Los datos del primer trimestre del año 2019 y posteriores son provisionales."
$ws.Range("E59").Value = "Los datos del primer trimestre del año 2019 y posteriores son provisionales."

# --- Row 16: add a note to E16 (was empty) ---
$ws.Range("E16").Value = "Para una mejor interpretación es el dato del indicador y no su tasa de variación anual"

# --- Row 74: update the "Tasa de apertura" note wording ---
$ws.Range("E74").Value = "Tasa de apertura= (saldo comercial/PIB)*100.  Para una mejor interpretación es el dato del indicador y no su tasa de variación anual"

# --- Row 54: E54 note changes and picks up the border/center style used in column F of that row ---
$ws.Range("F54").Copy()
$ws.Range("E54").PasteSpecial(-4122)
$ws.Range("E54").Value = "Los datos de mayo del año 2019 y posteriores son provisionales"

# --- Column widths for D, E, F (notes columns got narrower / lost "best fit") ---
$ws.Columns.Item(4).ColumnWidth = 29.140625
$ws.Columns.Item(5).ColumnWidth = 37.140625
$ws.Columns.Item(6).ColumnWidth = 31.5703125

# --- Scroll / selection state: view now sits lower on the sheet, active cell is E74 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 61
$win.ScrollColumn = 1
$ws.Range("E74").Select()
